# Fix the typo "ocument" -> "document" in the intro paragraph.
# The original text lives in a single run; Word's insertion-point edit
# splits it into three runs: "Dit ", the inserted "d", and the remainder
# "ocument vat ...". We reproduce that by locating the paragraph, finding
# the character offset right before "ocument", inserting the missing "d",
# and then nudging formatting on the newly-created sub-ranges so Word
# keeps them as separate runs (matching the split seen in the diff)
# instead of re-merging them into one run.

$d = $word.ActiveDocument

# Locate the paragraph that contains the typo.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*ocument vat de stageperiode*") {
        $target = $p
        break
    }
}

$pStart = $target.Range.Start

# "Dit " is the first 4 characters of the paragraph; the typo "ocument"
# starts right after it.
$insertAt = $pStart + 4

# Insert the missing "d" right before "ocument".
$gap = $d.Range($insertAt, $insertAt)
$gap.InsertBefore("d")

# Split "Dit " into its own run.
$runDit = $d.Range($pStart, $insertAt)
$runDit.Bold = 1
$runDit.Bold = 0

# Split the inserted "d" into its own run, separate from the remainder.
$runD = $d.Range($insertAt, $insertAt + 1)
$runD.Bold = 1
$runD.Bold = 0
